$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.288.59'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.870.79'
$ws.Range("E3").Value = '  +0.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7079'
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.49'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07767'
$ws.Range("E8").Value = '  +1.27%  '

$ws.Range("E9").Value = '  -0.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.98'
$ws.Range("E10").Value = '  +1.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08390'
$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.866.54'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.242'
$ws.Range("E13").Value = '  +0.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7103'
$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("E15").Value = '  -0.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.298.05'
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.068'
$ws.Range("E17").Value = '  +2.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008180'
$ws.Range("E18").Value = '  +4.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.45'
$ws.Range("E19").Value = '  -1.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.21'
$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.110.72'
$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.743'
$ws.Range("E23").Value = '  -1.45%  '

$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("E25").Value = '  -0.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.25'
$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.003'
$ws.Range("E27").Value = '  +0.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.43'
$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.503'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.398'
$ws.Range("E30").Value = '  -0.09%  '

$ws.Range("E31").Value = '  -1.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.296'
$ws.Range("E32").Value = '  +1.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05327'
$ws.Range("E33").Value = '  +3.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.938'
$ws.Range("E34").Value = '  +1.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.176'
$ws.Range("E35").Value = '  +0.75%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7435'
$ws.Range("E36").Value = '  -7.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.697'
$ws.Range("E37").Value = '  +0.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01867'
$ws.Range("E38").Value = '  +0.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.231.49'
$ws.Range("E39").Value = '  +6.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.725'
$ws.Range("E40").Value = '  +0.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.556'
$ws.Range("E41").Value = '  +3.88%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8859'
$ws.Range("E42").Value = '  -1.28%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '109.37'
$ws.Range("E43").Value = '  +6.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.33'
$ws.Range("E44").Value = '  -1.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9995'
$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.015.55'
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5190'
$ws.Range("E47").Value = '  +0.17%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.793'
$ws.Range("E48").Value = '  +0.69%  '

$ws.Range("E49").Value = '  +2.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.415'
$ws.Range("E50").Value = '  +0.78%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4309'
$ws.Range("E51").Value = '  +0.28%  '
